# Generate Report for Handback
# ------------------------------------------------------------------
# The localization CI job finished a handback round-trip for both
# target languages (zh-cn, de-de): the handoff source file is now the
# "latest target file" too (nothing left to localize -> in sync with
# en-US), and a fresh handback xliff + timestamp was produced for each
# language. Mirror that onto the Overview + per-language sheets.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$mdFile   = "74ea9e28-b2f5-4609-a8c7-113c66e10282.md"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/oltest/blob/01df3f3af4482f3c0538878be00a80e79ea418dc/e2e/74ea9e28-b2f5-4609-a8c7-113c66e10282.md"
$newStatus = "Handed back: in sync with en-US"

# -------------------- Overview sheet --------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus

# -------------------- zh-cn sheet --------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("I2").Value = $mdFile
$zh.Range("J2").Value = "74ea9e28-b2f5-4609-a8c7-113c66e10282.03ea18f8c25e25beffb2065396c6f3a6e7f93e22.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-13 05:12:18"

$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, "", "", $mdFile)
$zh.Range("I2").Font.Color = 6610909
$zh.Range("I2").Font.Underline = $true

# -------------------- de-de sheet --------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("I2").Value = $mdFile
$de.Range("J2").Value = "74ea9e28-b2f5-4609-a8c7-113c66e10282.03ea18f8c25e25beffb2065396c6f3a6e7f93e22.de-de.xlf"
$de.Range("K2").Value = "2016-08-13 05:12:28"

$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, "", "", $mdFile)
$de.Range("I2").Font.Color = 6610909
$de.Range("I2").Font.Underline = $true

# -------------------- column widths --------------------
# Widen the status columns (longer "Handed back: in sync with en-US"
# text) and the newly-populated Latest Target File / Latest Handback
# File columns on both language sheets.
$ovw.Columns.Item(5).ColumnWidth = 29.166666666666664
$ovw.Columns.Item(6).ColumnWidth = 29.166666666666664

$zh.Columns.Item(3).ColumnWidth  = 29.166666666666664
$zh.Columns.Item(9).ColumnWidth  = 39.16666666666667
$zh.Columns.Item(10).ColumnWidth = 39.16666666666667

$de.Columns.Item(3).ColumnWidth  = 29.166666666666664
$de.Columns.Item(9).ColumnWidth  = 39.16666666666667
$de.Columns.Item(10).ColumnWidth = 39.16666666666667
